$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'61.856.90"
$ws.Range('E2').Value = "'  -1.31%  "
$ws.Range('D3').Value = "'2.910.43"
$ws.Range('E3').Value = "'  -2.05%  "
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = "'  +0.16%  "
$ws.Range('D5').Value = "'587.15"
$ws.Range('E5').Value = "'  -1.38%  "
$ws.Range('D6').Value = "'145.63"
$ws.Range('E6').Value = "'  +0.89%  "
$ws.Range('E7').Value = "'  +0.04%  "
$ws.Range('D8').Value = "'0.505"
$ws.Range('E8').Value = "'  +0.50%  "
$ws.Range('D9').Value = "'2.908.67"
$ws.Range('E9').Value = "'  -2.02%  "
$ws.Range('D10').Value = "'7.04"
$ws.Range('E10').Value = "'  -3.02%  "
$ws.Range('D11').Value = "'0.151"
$ws.Range('E11').Value = "'  +4.31%  "
$ws.Range('D12').Value = "'0.437"
$ws.Range('E12').Value = "'  -1.89%  "
$ws.Range('D13').Value = "'0.0000239"
$ws.Range('E13').Value = "'  +1.08%  "
$ws.Range('D14').Value = "'32.89"
$ws.Range('E14').Value = "'  -2.07%  "
$ws.Range('E15').Value = "'  -1.50%  "
$ws.Range('D16').Value = "'3.394.37"
$ws.Range('E16').Value = "'  -1.90%  "
$ws.Range('D17').Value = "'61.900.44"
$ws.Range('E17').Value = "'  -0.90%  "
$ws.Range('D18').Value = "'6.60"
$ws.Range('E18').Value = "'  -1.90%  "
$ws.Range('D19').Value = "'2.910.74"
$ws.Range('E19').Value = "'  -1.96%  "
$ws.Range('D20').Value = "'434.42"
$ws.Range('E20').Value = "'  -1.78%  "
$ws.Range('D21').Value = "'13.45"
$ws.Range('E21').Value = "'  -1.01%  "
$ws.Range('D22').Value = "'0.658"
$ws.Range('E22').Value = "'  -2.65%  "
$ws.Range('D23').Value = "'6.92"
$ws.Range('E23').Value = "'  -3.65%  "
$ws.Range('D24').Value = "'81.02"
$ws.Range('E24').Value = "'  -1.10%  "
$ws.Range('B25').Value = "'InternetComputer(DFINITY)"
$ws.Range('C25').Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range('D25').Value = "'11.81"
$ws.Range('E25').Value = "'  -1.77%  "
$ws.Range('B26').Value = "'RenderToken"
$ws.Range('C26').Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range('D26').Value = "'10.16"
$ws.Range('E26').Value = "'  -6.41%  "
$ws.Range('B27').Value = "'Dai"
$ws.Range('C27').Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range('D27').Value = "'1.00"
$ws.Range('E27').Value = "'  -0.01%  "
$ws.Range('B28').Value = "'Fetch.AI"
$ws.Range('C28').Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range('D28').Value = "'2.06"
$ws.Range('E28').Value = "'  -4.38%  "
$ws.Range('D29').Value = "'0.0000108"
$ws.Range('E29').Value = "'  +22.90%  "
$ws.Range('D30').Value = "'7.21"
$ws.Range('E30').Value = "'  +2.74%  "
$ws.Range('D31').Value = "'2.56"
$ws.Range('E31').Value = "'  -1.82%  "
$ws.Range('D32').Value = "'2.09"
$ws.Range('E32').Value = "'  -1.10%  "
$ws.Range('E33').Value = "'  +1.41%  "
$ws.Range('B34').Value = "'FirstDigitalUSD"
$ws.Range('C34').Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range('D34').Value = "'1.00"
$ws.Range('E34').Value = "'  -0.05%  "
$ws.Range('B35').Value = "'EthereumClassic"
$ws.Range('C35').Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range('D35').Value = "'25.91"
$ws.Range('E35').Value = "'  -2.34%  "
$ws.Range('D36').Value = "'0.973"
$ws.Range('E36').Value = "'  -2.40%  "
$ws.Range('D37').Value = "'3.05"
$ws.Range('E37').Value = "'  +2.88%  "
$ws.Range('D38').Value = "'5.50"
$ws.Range('E38').Value = "'  -2.41%  "
$ws.Range('D39').Value = "'49.06"
$ws.Range('E39').Value = "'  -1.54%  "
$ws.Range('E40').Value = "'  -0.72%  "
$ws.Range('D41').Value = "'8.36"
$ws.Range('E41').Value = "'  -3.23%  "
$ws.Range('E42').Value = "'  -3.75%  "
$ws.Range('D43').Value = "'0.272"
$ws.Range('E43').Value = "'  -3.40%  "
$ws.Range('D44').Value = "'38.65"
$ws.Range('E44').Value = "'  -0.92%  "
$ws.Range('D45').Value = "'2.699.16"
$ws.Range('E45').Value = "'  +0.07%  "
$ws.Range('D46').Value = "'134.41"
$ws.Range('E46').Value = "'  +0.44%  "
$ws.Range('D47').Value = "'0.0338"
$ws.Range('E47').Value = "'  -0.97%  "
$ws.Range('D48').Value = "'343.60"
$ws.Range('E48').Value = "'  -6.89%  "
$ws.Range('E50').Value = "'  -1.15%  "
$ws.Range('D51').Value = "'22.28"
$ws.Range('E51').Value = "'  -4.12%  "
